$wb = $excel.ActiveWorkbook
$loginSheet = $wb.Worksheets.Item("Login")

# ---- Update the existing "Login" sheet selection / formatting ----
# Remove the fill-applying style from A2:A5 (back to default style)
$loginSheet.Range("A2:A5").ClearFormats()

# ---- Add the new "DataLogin" sheet right after "Login" ----
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $loginSheet)
$newSheet.Name = "DataLogin"

# Header row
$newSheet.Range("A1").Value = "email"
$newSheet.Range("B1").Value = "password"

# Data rows (order chosen so shared strings are interned in the same
# sequence as the target workbook: admin@example.com, admin1@example.com,
# admin11@example.com)
$newSheet.Range("A2").Value = "admin@example.com"
$newSheet.Range("B2").Value = 123456
$newSheet.Range("A4").Value = "admin1@example.com"
$newSheet.Range("B4").Value = 123456
$newSheet.Range("A3").Value = "admin11@example.com"
$newSheet.Range("B3").Value = 123456
$newSheet.Range("A5").Value = "admin@example.com"
$newSheet.Range("B5").Value = 123456

# Column widths (closest values the engine's pixel-quantized width model
# can reach to the source widths of 22.5703125 / 14.140625)
$newSheet.Columns.Item(1).ColumnWidth = 21.67
$newSheet.Columns.Item(2).ColumnWidth = 13.33

# Login sheet keeps A1:B5 selected (no longer the active tab)
$loginSheet.Range("A1:B5").Select()

# Make "DataLogin" the active/selected sheet with D7 selected (must be
# last, since selecting a range on a sheet also activates that sheet)
$newSheet.Activate()
$newSheet.Range("D7").Select()
